$wb = $excel.ActiveWorkbook

# --- Update last-updated timestamp on Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:58 PM"

# --- Update Stock List sheet: a new row (CAPTRU-RE1) is inserted at row 2,
#     pushing existing rows 2-75 down to rows 3-76, and the former last row
#     (row 76, TRAVELFOOD) drops off the bottom of the sheet. ---
$ws = $wb.Worksheets.Item("Stock List")

for ($r = 76; $r -ge 3; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($src, 2).Value()
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($src, 3).Value()
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value()
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($src, 5).Value()
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($src, 8).Value()
}

# New top entry
$ws.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value = 5.67
$ws.Cells.Item(2, 5).Value = -11.9565
$ws.Cells.Item(2, 8).Value = 0
